$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 51 (the "ULTIMATE MILLIONS ($30)" row with the anomalous 2019-02-09 scrape date).
# This shifts all rows below it up by one.
$ws.Rows.Item(51).Delete()
